$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 26352
$ws.Range("E2").Value = 510675236770
$ws.Range("F2").Value = 9477488509
$ws.Range("G2").Value = -3.12171
$ws.Range("D3").Value = 1800.02
$ws.Range("E3").Value = 216363342154
$ws.Range("F3").Value = 8001002124
$ws.Range("G3").Value = -2.88059
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 83028796839
$ws.Range("F4").Value = 23491717970
$ws.Range("G4").Value = 0.01933
$ws.Range("D5").Value = 306.09
$ws.Range("E5").Value = 48304263379
$ws.Range("F5").Value = 566911502
$ws.Range("G5").Value = -2.33032
$ws.Range("D6").Value = 0.999692
$ws.Range("E6").Value = 29127165251
$ws.Range("F6").Value = 3253358607
$ws.Range("G6").Value = 0.0008899999999999999
$ws.Range("D7").Value = 0.455152
$ws.Range("E7").Value = 23601672587
$ws.Range("F7").Value = 1118517791
$ws.Range("G7").Value = -2.22397
$ws.Range("D8").Value = 0.364581
$ws.Range("E8").Value = 12772728730
$ws.Range("F8").Value = 189263165
$ws.Range("G8").Value = -1.61656
$ws.Range("D9").Value = 1797.72
$ws.Range("E9").Value = 11992041666
$ws.Range("F9").Value = 11319492
$ws.Range("G9").Value = -2.85186
$ws.Range("D10").Value = 0.070877
$ws.Range("E10").Value = 9886542067
$ws.Range("F10").Value = 318174237
$ws.Range("G10").Value = -2.80056
$ws.Range("D11").Value = 0.874905
$ws.Range("E11").Value = 8117055439
$ws.Range("F11").Value = 239341729
$ws.Range("G11").Value = -1.77351
$ws.Range("E12").Value = 7619938877
$ws.Range("F12").Value = 321103279
$ws.Range("G12").Value = -4.11632
$ws.Range("D13").Value = 0.077029
$ws.Range("E13").Value = 6958382596
$ws.Range("F13").Value = 370834050
$ws.Range("G13").Value = -1.9949
$ws.Range("D14").Value = 5.26
$ws.Range("E14").Value = 6497062513
$ws.Range("F14").Value = 104965747
$ws.Range("G14").Value = -2.54622
$ws.Range("D15").Value = 86.02
$ws.Range("E15").Value = 6278838365
$ws.Range("F15").Value = 845772160
$ws.Range("G15").Value = -5.82987
$ws.Range("D16").Value = 0.999597
$ws.Range("E16").Value = 5333590561
$ws.Range("F16").Value = 2178215477
$ws.Range("G16").Value = -0.02084
$ws.Range("D17").Value = 0.000008549999999999999
$ws.Range("E17").Value = 5031393684
$ws.Range("F17").Value = 128059992
$ws.Range("G17").Value = -4.18783
$ws.Range("D18").Value = 14.13
$ws.Range("E18").Value = 4726914375
$ws.Range("F18").Value = 139711690
$ws.Range("G18").Value = -3.91398
$ws.Range("D19").Value = 0.999606
$ws.Range("E19").Value = 4627028447
$ws.Range("F19").Value = 120165964
$ws.Range("G19").Value = -0.04279
$ws.Range("D20").Value = 26402
$ws.Range("E20").Value = 4118586111
$ws.Range("F20").Value = 78160924
$ws.Range("G20").Value = -3.0684
$ws.Range("D21").Value = 4.98
$ws.Range("E21").Value = 3749850755
$ws.Range("F21").Value = 46532617
$ws.Range("G21").Value = -2.02506
$ws.Range("E22").Value = 3294013234
$ws.Range("F22").Value = 223644
$ws.Range("G22").Value = 0.8761100000000001
$ws.Range("D23").Value = 6.34
$ws.Range("E23").Value = 3274838315
$ws.Range("F23").Value = 155026690
$ws.Range("G23").Value = -2.66607
$ws.Range("D24").Value = 10.41
$ws.Range("E24").Value = 3044569189
$ws.Range("F24").Value = 77186250
$ws.Range("G24").Value = -1.07175
$ws.Range("D25").Value = 1.97
$ws.Range("E25").Value = 2896790669
$ws.Range("F25").Value = 12202379
$ws.Range("G25").Value = -3.96794
$ws.Range("D26").Value = 150.26
$ws.Range("E26").Value = 2727306374
$ws.Range("F26").Value = 64984950
$ws.Range("G26").Value = -0.6789500000000001
$ws.Range("D27").Value = 44.71
$ws.Range("E27").Value = 2683746185
$ws.Range("F27").Value = 7564958
$ws.Range("G27").Value = -2.76316
$ws.Range("D28").Value = 17.78
$ws.Range("E28").Value = 2508960103
$ws.Range("F28").Value = 76149025
$ws.Range("G28").Value = -3.23918
$ws.Range("D29").Value = 0.086622
$ws.Range("E29").Value = 2319106249
$ws.Range("F29").Value = 37187637
$ws.Range("G29").Value = -1.74851
$ws.Range("D30").Value = 112.31
$ws.Range("E30").Value = 2177770417
$ws.Range("F30").Value = 65039739
$ws.Range("G30").Value = -2.81503
$ws.Range("E31").Value = 2099920815
$ws.Range("F31").Value = 27508249
$ws.Range("G31").Value = -4.29772
$ws.Range("D32").Value = 0.999519
$ws.Range("E32").Value = 2040785907
$ws.Range("F32").Value = 231161664
$ws.Range("G32").Value = -0.02085
$ws.Range("E33").Value = 1906259810
$ws.Range("F33").Value = 110355388
$ws.Range("G33").Value = -1.02616
$ws.Range("E34").Value = 1799041584
$ws.Range("F34").Value = 63778715
$ws.Range("G34").Value = 0.29092
$ws.Range("D35").Value = 0.050997
$ws.Range("E35").Value = 1604273142
$ws.Range("F35").Value = 20052660
$ws.Range("G35").Value = -2.79466
$ws.Range("D36").Value = 8.09
$ws.Range("E36").Value = 1600806471
$ws.Range("F36").Value = 67392754
$ws.Range("G36").Value = -4.75764
$ws.Range("D37").Value = 0.059916
$ws.Range("E37").Value = 1513456266
$ws.Range("F37").Value = 7728402
$ws.Range("G37").Value = -3.59154
$ws.Range("D38").Value = 100.57
$ws.Range("E38").Value = 1462866211
$ws.Range("F38").Value = 15991436
$ws.Range("G38").Value = -2.12653
$ws.Range("D39").Value = 1.58
$ws.Range("E39").Value = 1432699621
$ws.Range("F39").Value = 61574852
$ws.Range("G39").Value = -3.9758
$ws.Range("E40").Value = 1415630407
$ws.Range("F40").Value = 191200833
$ws.Range("G40").Value = -4.5947
$ws.Range("D41").Value = 0.01944057
$ws.Range("E41").Value = 1412961718
$ws.Range("F41").Value = 41235137
$ws.Range("G41").Value = 0.13137
$ws.Range("E42").Value = 1228376218
$ws.Range("F42").Value = 73239921
$ws.Range("G42").Value = -4.79747
$ws.Range("D43").Value = 0.155338
$ws.Range("E43").Value = 1124448097
$ws.Range("F43").Value = 48637583
$ws.Range("G43").Value = -4.94942
$ws.Range("D44").Value = 0.116614
$ws.Range("E44").Value = 1047130585
$ws.Range("F44").Value = 38885460
$ws.Range("G44").Value = -4.6736
$ws.Range("D45").Value = 0.090447
$ws.Range("E45").Value = 1035950281
$ws.Range("F45").Value = 42037
$ws.Range("G45").Value = -0.07207
$ws.Range("D46").Value = 0.999682
$ws.Range("E46").Value = 1018676333
$ws.Range("F46").Value = 22087066
$ws.Range("G46").Value = 0.03234
$ws.Range("D47").Value = 0.998424
$ws.Range("E47").Value = 1002536398
$ws.Range("F47").Value = 10677254
$ws.Range("G47").Value = -0.14811
$ws.Range("D48").Value = 0.99992
$ws.Range("E48").Value = 999871469
$ws.Range("F48").Value = 5649721
$ws.Range("G48").Value = 0.02404
$ws.Range("E49").Value = 973789506
$ws.Range("F49").Value = 201810285
$ws.Range("G49").Value = -0.68514
$ws.Range("D50").Value = 0.840215
$ws.Range("E50").Value = 928634104
$ws.Range("F50").Value = 97822860
$ws.Range("G50").Value = -2.53788
$ws.Range("D51").Value = 0.331162
$ws.Range("E51").Value = 924344807
$ws.Range("F51").Value = 273297270
$ws.Range("G51").Value = -9.07572
